# tests/models/features/constants/input_data/input_data.xlsx
# Rename the "f_Name" header/column label to "flows_Name" and move the
# selected cell from F12 to E11 (matches the new selection saved in the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Q")

$ws.Range("B1").Value = "flows_Name"

$ws.Range("E11").Select()
